$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2  = @(0.003994804209775715, 0.002777888934908601, 0.1575252929769615, 0.496779210170732, 0.6610771962923778)
    3  = @(0.0008583669626518464, 0.3127903958511391, 0.1575252929769615, 0.496779210170732, 0.9679532659614845)
    4  = @(3.230985683306322, 1.667794583268128, 0.8054896365839992, 0.496779210170732, 6.201049113329182)
    5  = @(3.230985683306322, 1.667794583268128, 3.900430680208489, 0.496779210170732, 9.295990156953671)
    6  = @(1.459612070389937, 1.667794583268128, 0.8054896365839992, 0.496779210170732, 4.429675500412797)
    7  = @(1.459612070389937, 1.667794583268128, 0.8054896365839992, 0.496779210170732, 4.429675500412797)
    8  = @(3.230985683306322, 1.667794583268128, 0.8054896365839992, 0.496779210170732, 6.201049113329182)
    9  = @(3.230985683306322, 1.667794583268128, 0.8054896365839992, 0.496779210170732, 6.201049113329182)
    10 = @(0.04763786555579896, 0.3127903958511391, 3.900430680208489, 8.660232485948974, 12.9210914275644)
    11 = @(3.230985683306322, 10.29869402782916, 3.900430680208489, 8.660232485948974, 26.09034287729295)
    12 = @(3.230985683306322, 1.667794583268128, 0.1575252929769615, 0.496779210170732, 5.553084769722144)
    13 = @(0.3048080303191223, 0.3127903958511391, 0.1575252929769615, 0.496779210170732, 1.271902929317955)
    14 = @(1.459612070389937, 1.667794583268128, 0.1575252929769615, 0.496779210170732, 3.781711156805759)
    15 = @(3.230985683306322, 1.667794583268128, 0.8054896365839992, 0.496779210170732, 6.201049113329182)
    16 = @(1.459612070389937, 1.667794583268128, 0.8054896365839992, 0.496779210170732, 4.429675500412797)
    17 = @(1.459612070389937, 1.667794583268128, 26.21740644021617, 0.496779210170732, 29.84159230404497)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 2).Value = $vals[0]
    $ws.Cells.Item($row, 3).Value = $vals[1]
    $ws.Cells.Item($row, 4).Value = $vals[2]
    $ws.Cells.Item($row, 5).Value = $vals[3]
    $ws.Cells.Item($row, 7).Value = $vals[4]
}
